# DeliveryChanges.xlsx update
# - removes the vertical-center alignment that had been applied to the
#   data range (A1:D22)
# - clears the lingering E7 selection, leaving the sheet at the default A1
# - appends a new trailing row (23) with two numeric values (123, 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "vertical center" alignment from every populated cell while
# leaving the existing number formats (dates in column A, #,##0 in some
# column D cells) untouched.
$ws.Range("A1:D22").VerticalAlignment = -4107   # xlVAlignBottom (Excel default)

# New data row appended below the existing table.
$ws.Range("A23").Value = 123
$ws.Range("B23").Value = 4

# Reset the active selection to A1 (the previous file had a stray E7
# selection saved in it).
$ws.Range("A1").Select() | Out-Null
